# Tabla de simbolos - reassign the "Asignacion" / "Operacion Aritmetica" /
# "Delimitar" token classifications for rows 14 (the "=" row) and 17-18
# (the "(" / ")" rows), matching the corrected classification in the
# author's re-upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tabla de simbolos")

# Row 14 ("=" / igual): was Tipo_2=Operacion Aritmetica, Tipo_3=Asignacion.
# Now it is simply Tipo_2=Asignacion, with Tipo_3 cleared.
$ws.Range("E14").Value = "Asignacion"
$ws.Range("F14").ClearContents()

# Row 17 ("(" / parentIzq): was Tipo_2=" Delimitar" only.
# Now Tipo_2=Operacion Aritmetica, Tipo_3=" Delimitar".
$ws.Range("E17").Value = "Operacion Aritmetica"
$ws.Range("F17").Value = " Delimitar"

# Row 18 (")" / parentDer): same reshuffle as row 17.
$ws.Range("E18").Value = "Operacion Aritmetica"
$ws.Range("F18").Value = " Delimitar"

# Reflect the author's on-screen selection/scroll state at save time.
$ws.Activate()
try {
    $win = $excel.ActiveWindow
    if ($win) {
        $win.ScrollRow = 9
        $win.ScrollColumn = 1
    }
} catch {
    # Window-scroll position is cosmetic; ignore if unsupported.
}
$ws.Range("F15:F16").Select() | Out-Null
